$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" ---------------
# "Overview" sheet: columns E (zh-cn) and F (de-de), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# "zh-cn" sheet: column C (Status), rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# "de-de" sheet: column C (Status), rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Narrow the "Status"-holding columns ---------------------------------
# Original stored width 17.2159881591797 -> target stored width 13.4101845877511.
# The engine quantizes ColumnWidth to 1/6-character (pixel) steps, so a
# ColumnWidth of 12.5 is the closest achievable input to the target width.
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
